# Process 2 (Main) Updated - Archive sheet URL list maintenance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Archive")

# Replace the Zomato article URL (now tracked on the "Arcive" sheet) with the
# new OnePlus article URL in the Archive list.
$ws.Range("A4").Value2 = "https://gadgets.ndtv.com/mobiles/news/oneplus-8-pro-sale-india-today-june-29-price-rs-54989-specifications-amazon-offer-2253327"

# Append newly archived article URLs.
$ws.Range("A8").Value2  = "https://www.entrepreneur.com/article/347406"
$ws.Range("A9").Value2  = "https://www.bbc.co.uk/news/entertainment-arts-53190585"
$ws.Range("A10").Value2 = "https://news.sky.com/story/molly-conlin-former-eastenders-star-held-at-knifepoint-by-robbers-who-broke-into-her-home-12017105"

# Leave the active selection on the Archive sheet where the user finished working.
$ws.Activate()
$ws.Range("F8").Select()
